# SAM TODO.xlsx update
#  - Insert 3 new rows (88-90) on the "To Do" sheet with new TODO items
#  - Mark the "Direct steam power tower" item (row 53) as Done
#  - Update the saved selection/active cell on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# --- Insert three new rows above the current row 88 (shifts 88-128 -> 91-131) ---
$ws.Rows("88:90").Insert()

# --- Populate the three new rows with the new TODO items ---
$ws.Range("A88").Value = "Not done"
$ws.Range("B88").Value = "SDK build"
$ws.Range("C88").Value = "Steve"

$ws.Range("A89").Value = "Not done"
$ws.Range("B89").Value = "CEC Inverter updates"
$ws.Range("C89").Value = "Steve"

$ws.Range("A90").Value = "Not done"
$ws.Range("B90").Value = "PBNS update for dispatch factors"
$ws.Range("C90").Value = "Steve"

# --- Mark "Direct steam power tower" (row 53) as Done ---
$ws.Range("A53").Value = "Done"

# --- Update the sheet's saved selection ---
$ws.Range("C90").Select()
